$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts existing rows 13-16 down to 14-17).
$ws.Rows.Item(13).Insert() | Out-Null

# Copy the formatting (bold/border/center style) of the column-A cell above
# down into the new row's A13 cell, so it keeps style index 1 like the rest
# of the table's first column.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

# New roster entry for the upcoming matchday - only the player name is known
# so far, the score columns (C:AI) are intentionally left blank.
$ws.Range("B13").Value = "Муратов Игорь"

# Refresh the remembered "Data > Sort" range/condition so it covers the
# table's new extent (one additional data row).
$sort = $ws.Sort
$sort.SortFields.Clear() | Out-Null
$sort.SortFields.Add($ws.Range("B1:B17")) | Out-Null
$sort.SetRange($ws.Range("A2:Q17")) | Out-Null
$sort.Header = 0
$sort.Apply() | Out-Null

# Second new roster entry, appended after the table as a fresh row (18),
# again with only the name filled in.
$ws.Range("B18").Value = "Шевчук Антон"

# Match the cursor position left behind by the edit.
$ws.Range("H13").Select() | Out-Null
